# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> bound to the slide master ("Integral" theme)
#   ppt/theme/theme2.xml -> bound to the notes master   ("Office Theme")
# The target revision swaps which theme is which: the slide master's
# theme becomes the stock "Office Theme" color palette (and the notes
# master's theme becomes "Integral"). The font scheme and format scheme
# (fills/lines/effects) are already byte-identical between the two theme
# parts, so the only real content difference is the 12-slot color scheme
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
#
# Apply that by rewriting the slide (master) theme's color scheme, slot
# by slot, to the "Office Theme" palette via the ThemeColorScheme object
# that PowerPoint exposes off a Slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order is fixed: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6,
# 11 hlink, 12 folHlink. Values are standard VBA RGB() encodings
# (0x00BBGGRR, i.e. byte-reversed hex) of the target "Office Theme" hex
# colors (RGB() itself isn't available in this host, so the already
# computed decimal values are used directly).
$tcs.Item(1).RGB  = 0            # dk1      000000
$tcs.Item(2).RGB  = 16777215     # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388      # dk2      44546A
$tcs.Item(4).RGB  = 15132391     # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939     # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501      # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845     # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407        # accent4  FFC000
$tcs.Item(9).RGB  = 12874308     # accent5  4472C4
$tcs.Item(10).RGB = 4697456      # accent6  70AD47
$tcs.Item(11).RGB = 12673797     # hlink    0563C1
$tcs.Item(12).RGB = 7491477      # folHlink 954F72
